$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (12) and clean up the stray
# formatted cell that Excel's auto-extend creates in the merged F3:K3 row.
$ws.Range("L1").EntireColumn.Insert()
$ws.Range("L3").Clear()

# New header cell for the inserted column, matching the other row-8 headers.
$ws.Range("L8").Value = "DESCRIPCION APLICACIÓN SAE"

# Match the column's best-fit width for the new header text.
$ws.Columns("L").ColumnWidth = 28.42578125

# Update the active selection to match the post-edit workbook state.
$ws.Range("AG9").Select()
